# Clear the personal/record data that was entered into row 2 of the sheet
# (Name, Account no., Father name, Mother name, Address, Gender, DOB, Age,
# Phone, Aadhar, Security Pin, Date of creation) while leaving the header
# row (row 1) and the sheet's layout/formatting untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:L2").ClearContents()
